# Agrego nueva liga a RockData
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows (8 rows) appended to the fixtures table.
$rows = @(
    @{ Row=138; Fecha="2025-07-25"; Local="Recoleta";                    Visita="Concepción";                  GL=2; GV=3; FID=1348372; CL=5; CV=5; AL=4; AV=4; RL=1; RV=0; G1L=1; G1V=1; G2L=1; G2V=2; PL=54; PV=46; Res="V" },
    @{ Row=139; Fecha="2025-07-26"; Local="Union San Felipe";            Visita="Deportes Temuco";              GL=1; GV=0; FID=1348371; CL=5; CV=3; AL=3; AV=4; RL=0; RV=1; G1L=1; G1V=0; G2L=0; G2V=0; PL=57; PV=43; Res="L" },
    @{ Row=140; Fecha="2025-07-26"; Local="Deportes Santa Cruz";         Visita="San Luis";                     GL=2; GV=2; FID=1348373; CL=5; CV=5; AL=2; AV=3; RL=0; RV=0; G1L=1; G1V=0; G2L=1; G2V=2; PL=51; PV=49; Res="E" },
    @{ Row=141; Fecha="2025-07-26"; Local="Curico Unido";                Visita="Santiago Morning";             GL=0; GV=1; FID=1348374; CL=6; CV=2; AL=5; AV=4; RL=2; RV=1; G1L=0; G1V=1; G2L=0; G2V=0; PL=54; PV=46; Res="V" },
    @{ Row=142; Fecha="2025-07-26"; Local="Universidad de Concepcion";   Visita="Cobreloa";                     GL=0; GV=1; FID=1348376; CL=4; CV=5; AL=5; AV=5; RL=0; RV=0; G1L=0; G1V=1; G2L=0; G2V=0; PL=49; PV=51; Res="V" },
    @{ Row=143; Fecha="2025-07-27"; Local="Antofagasta";                 Visita="Santiago Wanderers";           GL=4; GV=1; FID=1348370; CL=2; CV=2; AL=4; AV=5; RL=1; RV=0; G1L=3; G1V=1; G2L=1; G2V=0; PL=56; PV=44; Res="L" },
    @{ Row=144; Fecha="2025-07-27"; Local="Rangers de Talca";            Visita="Deportes Copiapo";             GL=0; GV=0; FID=1348375; CL=4; CV=6; AL=1; AV=5; RL=0; RV=0; G1L=0; G1V=0; G2L=0; G2V=0; PL=56; PV=44; Res="E" },
    @{ Row=145; Fecha="2025-07-29"; Local="San Marcos de Arica";         Visita="Magallanes";                   GL=2; GV=2; FID=1348369; CL=3; CV=2; AL=3; AV=4; RL=1; RV=0; G1L=0; G1V=1; G2L=2; G2V=1; PL=50; PV=50; Res="E" }
)

foreach ($r in $rows) {
    $row = $r.Row
    # Column A holds a text date like "2025-07-26". Typing that literally into
    # Excel gets auto-recognised as a date serial, so instead build it with a
    # text formula and flatten the formula down to its literal string result
    # (below) -- this keeps the cell a plain shared-string, same as every
    # other "Fecha" cell in the sheet.
    $ws.Cells.Item($row, 1).Formula = '="' + $r.Fecha + '"'

    $ws.Cells.Item($row, 2).Value = $r.Local
    $ws.Cells.Item($row, 3).Value = $r.Visita
    $ws.Cells.Item($row, 4).Value = $r.GL
    $ws.Cells.Item($row, 5).Value = $r.GV
    $ws.Cells.Item($row, 6).Value = $r.FID
    $ws.Cells.Item($row, 7).Value = $r.CL
    $ws.Cells.Item($row, 8).Value = $r.CV
    $ws.Cells.Item($row, 9).Value = $r.AL
    $ws.Cells.Item($row, 10).Value = $r.AV
    $ws.Cells.Item($row, 11).Value = $r.RL
    $ws.Cells.Item($row, 12).Value = $r.RV
    $ws.Cells.Item($row, 13).Value = $r.G1L
    $ws.Cells.Item($row, 14).Value = $r.G1V
    $ws.Cells.Item($row, 15).Value = $r.G2L
    $ws.Cells.Item($row, 16).Value = $r.G2V
    $ws.Cells.Item($row, 17).Value = $r.PL
    $ws.Cells.Item($row, 18).Value = $r.PV
    $ws.Cells.Item($row, 19).Value = $r.Res
}

# Flatten the helper formulas in column A down to their literal text values
# (copy + paste-values), exactly like the rest of the "Fecha" column.
$dateRange = $ws.Range("A138:A145")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Column A now holds entries as wide as the existing dates -- refresh the
# best-fit width like Excel does automatically when new data is entered.
$ws.Columns("A").AutoFit()

# Reflect where the user ended up after entering the new fixtures.
$ws.Range("R145").Select()
